# Apply the "custom cutoffs" label update to Sheet1, column A.
# The values in column B and C are left untouched; only the row labels
# (shared-string text) in column A change to the new, more granular set
# of cutoff names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "cutoff"
$ws.Range("A2").Value = "_ipReads_cutoff"
$ws.Range("A3").Value = "_trimmedReads_cutoff"
$ws.Range("A4").Value = "_uniqAligned_cutoff"
$ws.Range("A5").Value = "_exonMapping_cutoff"
$ws.Range("A6").Value = "_riboScatter_cutoff"
$ws.Range("A7").Value = "_violin_cutoff_overrep_untrimmed"
$ws.Range("A9").Value = "_violin_cutoff_overrep_trimmed"
$ws.Range("A8").Value = "_violin_cutoff_adapter_untrimmed"
$ws.Range("A10").Value = "_violin_cutoff_adapter_trimmed"
$ws.Range("A11").Value = "GeneBody_Coverage"
$ws.Range("A12").Value = "Dist_of_gene_expression"

# Move the cell cursor / selection to A14, mirroring the session-state
# change captured in the sheetView (new <selection activeCell="A14".../>).
$ws.Range("A14").Select()
